# Milestone 5 - final updates
# Applies the set of text edits described by the commit diff.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARN: replace failed for: $find"
    }
}

# 1. "taking seasonality into account" -> "taking seasonality in to account"
Replace-Text "taking seasonality into account" "taking seasonality in to account"

# 2. "tune the models hyper-parameters." -> "tune the model's hyper-parameters."
Replace-Text "tune the models hyper-parameters." "tune the model’s hyper-parameters."

# 3. Rewrite of the "We'll conclude..." sentence.
Replace-Text "with a visualization of the feature importance (from the random forest model assuming we have a good fit to the data) to highlight which factors contribute to price. We will also include a visualization demonstrating goodness of fit (e.g. an ROC curve). For our website we also plan to include a static or somewhat staged demonstration of how our model could be used in the Airbnb UI." "with a visualization of feature importance (from the random forest model, assuming we have a reasonable fit to the data) to highlight which factors contribute to price. We’ll also include a visualization demonstrating goodness of fit (e.g. an ROC curve, suggestions welcome!). For our website, we also plan to include a static or somewhat staged demonstration of how our model could be used in the Airbnb UI."

# 4. "make predictions on price that beat" -> "make predictions on the price category that beat"
Replace-Text "determine whether or not we can make predictions on price that beat our baselines." "determine whether or not we can make predictions on the price category that beat our baselines."

# 5. "description are the most important" -> "description are most important"
Replace-Text "description are the most important when it comes to predicting prices." "description are most important when it comes to predicting prices."

# 6. Move the "_GoBack" bookmark from the end of the "We'll conclude..." paragraph
#    to just before "most important" in the "Since our goal is to suggest..." paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$findRange = $d.Content
$found = $findRange.Find.Execute("description are ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Range($findRange.End, $findRange.End)
    $d.Bookmarks.Add("_GoBack", $target)
} else {
    Write-Output "WARN: could not find bookmark insertion point"
}

Write-Output "Text replacements complete."
